$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reshuffles the data rows (sheet rows 3-14) of the "artfynd" export:
# every row's full content moves to a different row number (rows 1-2 stay put).
# Mapping is keyed by destination row -> source row (both in the original sheet).
$mapping = @{
    3  = 7
    4  = 12
    5  = 10
    6  = 14
    7  = 9
    8  = 6
    9  = 3
    10 = 11
    11 = 8
    12 = 13
    13 = 4
    14 = 5
}

$firstRow = 3
$lastRow  = 14
$lastCol  = 51   # column AY

# Startdatum/Starttid/Slutdatum/Sluttid (Y:AB) are stored as plain text that
# looks like dates/times (e.g. "2023-08-14", "00:00"). Force the range to
# Text format first so re-assigning the captured values below doesn't let
# Excel reinterpret them as real date/time serials.
$ws.Range("Y${firstRow}:AB${lastRow}").NumberFormat = "@"

# Snapshot all the source rows before overwriting anything.
$data = $ws.Range("A${firstRow}:AY${lastRow}").Value2

$rowCount = ($lastRow - $firstRow) + 1
$new = New-Object 'object[,]' $rowCount, $lastCol

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $destIdx = $destRow - $firstRow
    $srcIdx  = $srcRow - $firstRow
    for ($col = 0; $col -lt $lastCol; $col++) {
        $new[$destIdx, $col] = $data[$srcIdx + 1, $col + 1]
    }
}

$ws.Range("A${firstRow}:AY${lastRow}").Value2 = $new
